$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add two new rows of data for the two new certificate types
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Anagrafico di Unione Civile"

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "di Contratto di Convivenza"

# Update the active cell selection to match the saved workbook state
$ws.Range("D23").Select()
